$p = $ppt.ActivePresentation

# Slide 1: Set speaker notes text (LangChain4j overview)
$notes1 = @'
【スライド1: LangChain4jの概要】
＜イントロダクション - 2分＞
皆さん、こんにちは。本日は、JavaアプリケーションにAI機能を統合するための強力なフレームワーク、LangChain4jについてご紹介します。
＜現状の課題 - 1分＞
多くのJava開発者の方が、ChatGPTやClaudeなどのLLMをアプリケーションに組み込みたいと考えていますが、以下の課題に直面しています：
• Pythonベースのツールが多く、Javaエコシステムとの統合が困難
• 各AIプロバイダーのAPIが異なり、学習コストが高い
• エンタープライズ要件（セキュリティ、監査、スケーラビリティ）への対応が複雑
＜LangChain4jとは - 2分＞
LangChain4jは、これらの課題を解決するJava専用のLLMフレームワークです。
• PythonのLangChainの設計思想を継承しつつ、Java開発者のための最適化を実施
• 「Write once, run with any LLM」- 一度書けば、どのLLMでも動作
＜3つの主要特徴の説明 - 3分＞
1. 簡単な統合
   「既存のSpring BootアプリケーションにChatGPT機能を追加するのに、わずか10行のコードで実現できます」
   実例：カスタマーサポートチャットボット、コード生成ツール
2. チェーン構築
   「複数のAI処理を連結し、高度なワークフローを構築できます」
   実例：文書要約→感情分析→返信生成のような複合処理
3. Java最適化
   「StreamやOptional、アノテーションベースの設定など、Javaらしい実装が可能」
   Spring開発者なら違和感なく使い始められます
＜バージョン情報 - 30秒＞
• Java 8以上で動作（レガシーシステムでも導入可能）
• 活発な開発（月1-2回のリリース）
• Apache License 2.0（商用利用も安心）
＜想定される質問と回答＞
Q: 「Pythonのライブラリを使う方が良いのでは？」
A: 「既存のJavaインフラ、CI/CD、監視ツールをそのまま活用できる点が大きなメリットです。チーム全体がJavaに精通している場合、学習コストも最小限です。」
Q: 「パフォーマンスはどうですか？」
A: 「JVMの最適化により、特に大規模バッチ処理では優れたパフォーマンスを発揮します。」
'@
$p.Slides.Item(1).NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = $notes1

# Slide 2: Set speaker notes text (Components and architecture)
$notes2 = @'
【スライド2: コンポーネントとアーキテクチャ】
＜スライド切り替え - 30秒＞
「それでは、LangChain4jの技術的な構成について、詳しく見ていきましょう。」
＜主要コンポーネントの説明 - 5分＞
1. Language Models（言語モデル）- 1分30秒
   「LangChain4jの心臓部です。複数のAIプロバイダーを統一インターフェースで扱えます」
   対応プロバイダー：
   • OpenAI (GPT-4, GPT-3.5)
   • Anthropic (Claude)
   • Azure OpenAI Service
   • Google Vertex AI
   • Hugging Face
   • Ollama（ローカルLLM）
   実装例：「プロバイダーを変更する際も、設定ファイルを1行変更するだけです」
   コスト最適化：「開発環境ではGPT-3.5、本番ではGPT-4といった使い分けが簡単」
2. Embeddings & Vector Store - 1分30秒
   「RAG（Retrieval-Augmented Generation）を実現するための重要なコンポーネントです」
   使用シナリオ：
   • 社内文書検索システム
   • FAQチャットボット
   • コードベースの類似検索
   対応ベクトルDB：
   • インメモリ（開発用）
   • ChromaDB、Pinecone（クラウド）
   • Elasticsearch（既存インフラ活用）
   「10万件の文書でも、ミリ秒単位で関連情報を検索できます」
3. Document Processing - 1分
   「様々な形式のドキュメントをAIが理解できる形に変換します」
   対応フォーマット：
   • PDF（契約書、マニュアル）
   • Word、Excel（業務文書）
   • HTML（Webコンテンツ）
   • テキストファイル（ログ、設定ファイル）
   「PDFの表も正確に抽出し、構造を保持したまま処理できます」
4. Chains & Memory - 1分
   「複数の処理を連結し、文脈を保持した対話を実現します」
   メモリタイプ：
   • ConversationBufferMemory（全履歴保持）
   • ConversationSummaryMemory（要約保持）
   • ConversationWindowMemory（直近N件）
   「カスタマーサポートでは、過去の会話を参照しながら適切な回答が可能」
＜アーキテクチャ階層の説明 - 2分＞
「4層構造により、関心の分離と拡張性を実現しています」
• アプリケーション層
  「ビジネスロジックとUIを実装。REST API、WebSocketでの公開も容易」
• コア層
  「フレームワークの中核。プロンプトテンプレート、チェーン実行エンジンなど」
• プロバイダー層
  「各AIサービスとの接続。新しいプロバイダーの追加も簡単」
• インフラストラクチャー層
  「永続化、キャッシング、監視などの横断的関心事」
＜実装のポイント - 1分＞
「Spring Bootを使用している場合、@Configurationクラスでの設定だけで開始できます」
「既存のDIコンテナ、トランザクション管理、セキュリティ設定をそのまま活用可能」
＜想定される質問と回答＞
Q: 「ベクトルストアの選定基準は？」
A: 「データ量、レイテンシ要件、既存インフラを考慮。開発はインメモリ、本番はマネージドサービスを推奨」
Q: 「メモリ使用量が心配です」
A: 「ConversationWindowMemoryやSummaryMemoryで制御可能。Redisでの外部化も対応」
'@
$p.Slides.Item(2).NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = $notes2

# Slide 3: Set speaker notes text (Usage and benefits)
$notes3 = @'
【スライド3: 使用方法とメリット】
＜導入 - 30秒＞
「最後に、実際のコード例とビジネス価値についてご説明します。」
＜コード例の詳細解説 - 4分＞
1. セットアップ（1分）
   「Maven依存関係を追加するだけで始められます」
   依存関係の説明：
   • langchain4j-core: 基本機能
   • langchain4j-open-ai: OpenAI連携（必要に応じて選択）
   • langchain4j-spring-boot-starter: Spring Boot統合
   「既存プロジェクトへの影響は最小限。JARサイズも約2MB程度」
2. 初期化コード（1分30秒）
   ChatLanguageModel model = OpenAiChatModel.builder()
       .apiKey(System.getenv("OPENAI_API_KEY"))
       .modelName(OpenAiModelName.GPT_4)
       .temperature(0.7)
       .build();
   ポイント説明：
   • 「APIキーは環境変数から読み込み（セキュリティベストプラクティス）」
   • 「temperatureパラメータで創造性を調整（0.0=決定的、1.0=創造的）」
   • 「タイムアウト、リトライ、プロキシ設定も可能」
3. 実行例（1分30秒）
   String response = model.generate("JavaでFizzBuzzを実装してください");
   「たった1行で、自然言語での指示をコードに変換できます」
   高度な使用例：
   • ストリーミング: model.generateStream()でリアルタイム応答
   • 関数呼び出し: 外部APIとの連携
   • プロンプトテンプレート: 再利用可能な質問形式
＜ビジネスメリットの詳細 - 5分＞
1. 開発速度向上（1分30秒）
   「実際の導入事例では、開発工数を30-50%削減」
   具体例：
   • チャットボット開発: 3ヶ月→1ヶ月
   • ドキュメント検索システム: 2ヶ月→3週間
   • コードレビューツール: 6週間→2週間
   「ボイラープレートコードが不要。ビジネスロジックに集中できます」
2. プロバイダー非依存（1分30秒）
   「将来の技術変化にも柔軟に対応」
   切り替えシナリオ：
   • コスト削減: OpenAI→Claude（50%コスト削減の事例）
   • パフォーマンス: GPT-3.5→GPT-4（精度向上）
   • コンプライアンス: クラウド→オンプレミス
   「設定ファイルの変更だけで、コードの修正は不要」
3. エンタープライズ対応（2分）
   「大企業の厳格な要件にも対応」
   セキュリティ機能：
   • Spring Security統合（認証・認可）
   • APIキーの暗号化管理
   • プロンプトインジェクション対策
   監査・コンプライアンス：
   • 全リクエスト/レスポンスのロギング
   • GDPR対応（個人情報のマスキング）
   • SOC2準拠のための監査証跡
   スケーラビリティ：
   • 非同期処理対応
   • コネクションプーリング
   • 分散キャッシング（Redis統合）
＜活用事例 - 2分＞
「実際の導入企業での成功事例をご紹介します」
1. 金融機関A社
   「コンプライアンス文書の自動チェックシステム」
   • 処理時間: 90%削減
   • 精度: 人間のチェックと同等以上
2. Eコマース企業B社
   「商品説明の自動生成」
   • 月間10万商品の説明文を自動作成
   • SEOスコア: 平均30%向上
3. SaaS企業C社
   「カスタマーサポートの自動化」
   • 問い合わせの70%を自動解決
   • 顧客満足度: 15%向上
＜実装のベストプラクティス - 1分＞
• エラーハンドリング: Circuit Breakerパターンの実装
• コスト管理: トークン使用量の監視とアラート
• プロンプト管理: バージョン管理とA/Bテスト
• キャッシング: 同一質問への応答をキャッシュ
＜想定される質問と回答＞
Q: 「導入にあたっての最初のステップは？」
A: 「まず小さなPoCから始めることを推奨。社内FAQボットなど、リスクの低い領域から導入し、段階的に拡大」
Q: 「運用コストはどの程度？」
A: 「月間100万リクエストで約$500-1000。キャッシングとプロンプト最適化でさらに削減可能」
Q: 「社内データの学習は可能？」
A: 「RAGパターンで社内文書を参照可能。Fine-tuningも対応予定」
＜クロージング - 30秒＞
「LangChain4jは、Javaエコシステムに最適化された実用的なLLMフレームワークです。
ぜひ、小さなプロジェクトから始めて、その威力を実感してください。
ご質問があれば、お気軽にお聞きください。」
'@
$p.Slides.Item(3).NotesPage.Shapes.Item(2).TextFrame.TextRange.Text = $notes3

# Update gradient image alt-text/description identifiers on each slide's background image
$p.Slides.Item(1).Shapes.Item(3).AlternativeText = "/tmp/rasterized-gradient-0bd548ff.png"
$p.Slides.Item(2).Shapes.Item(21).AlternativeText = "/tmp/rasterized-gradient-cf4f7bce.png"
$p.Slides.Item(3).Shapes.Item(31).AlternativeText = "/tmp/rasterized-gradient-20f94e08.png"
